$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.035900592803955
$ws.Range("B1").Value = 2.420019626617432
$ws.Range("C1").Value = 5.133058547973633
$ws.Range("D1").Value = 2.31307315826416
$ws.Range("E1").Value = 1.328786730766296
